$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1024
$ws.Range("J70").Value = 1024
$ws.Range("L70").Value = 3072
$ws.Range("N70").Value = -3612
$ws.Range("H73").Value = 1024
$ws.Range("J73").Value = 1024
$ws.Range("L73").Value = 3072
$ws.Range("N73").Value = -4944
$ws.Range("H99").Value = 273.33334
$ws.Range("J99").Value = 499
$ws.Range("L99").Value = 1497
$ws.Range("N99").Value = -4493
$ws.Range("H101").Value = 16669804
$ws.Range("J101").Value = 475
$ws.Range("L101").Value = 1425
$ws.Range("N101").Value = -4669
$ws.Range("H116").Value = 10701
$ws.Range("I116").Value = 11176.25
$ws.Range("J116").Value = 8800
$ws.Range("K116").Value = 11176.25
$ws.Range("L116").Value = 8800
$ws.Range("M116").Value = -7734.25
$ws.Range("N116").Value = -15684
$ws.Range("H137").Value = 1484.5
$ws.Range("I137").Value = 1479
$ws.Range("K137").Value = 4437
$ws.Range("M137").Value = -1887
$ws.Range("H138").Value = 4624.8335
$ws.Range("I138").Value = 2119.75
$ws.Range("K138").Value = 6359.25
$ws.Range("M138").Value = -1219.25

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1598.4
$ws.Range("I2").Value = 1598.4
$ws.Range("K2").Value = 1598.4
$ws.Range("M2").Value = -1485.4
$ws.Range("H10").Value = 502974.75
$ws.Range("I10").Value = 502974.75
$ws.Range("K10").Value = 502974.75
$ws.Range("M10").Value = -502804.75
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("H74").Value = 1683
$ws.Range("I74").Value = 1913.6666
$ws.Range("J74").Value = 299
$ws.Range("K74").Value = 1913.6666
$ws.Range("L74").Value = 299
$ws.Range("M74").Value = -1039.6666
$ws.Range("N74").Value = -2047
$ws.Range("H77").Value = 1683
$ws.Range("I77").Value = 1913.6666
$ws.Range("J77").Value = 299
$ws.Range("K77").Value = 9568.333000000001
$ws.Range("L77").Value = 1495
$ws.Range("M77").Value = -5200.333000000001
$ws.Range("N77").Value = -10231
$ws.Range("H97").Value = 570.2308
$ws.Range("I97").Value = 590.4
$ws.Range("K97").Value = 590.4
$ws.Range("M97").Value = -94.39999999999998
$ws.Range("H110").Value = 3078.4
$ws.Range("I110").Value = 3078.4
$ws.Range("K110").Value = 3078.4
$ws.Range("M110").Value = -1033.4
$ws.Range("H116").Value = 1598.4
$ws.Range("I116").Value = 1598.4
$ws.Range("K116").Value = 1598.4
$ws.Range("M116").Value = 695.5999999999999
$ws.Range("H122").Value = 1532.3334
$ws.Range("I122").Value = 1298
$ws.Range("J122").Value = 1649.5
$ws.Range("K122").Value = 3894
$ws.Range("L122").Value = 4948.5
$ws.Range("M122").Value = -1444
$ws.Range("N122").Value = -9848.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1598.4
$ws.Range("I3").Value = 1598.4
$ws.Range("K3").Value = 1598.4
$ws.Range("M3").Value = -1484.4
$ws.Range("H20").Value = 1865
$ws.Range("I20").Value = 853.6667
$ws.Range("J20").Value = 3078.6
$ws.Range("K20").Value = 853.6667
$ws.Range("L20").Value = 3078.6
$ws.Range("M20").Value = -606.6667
$ws.Range("N20").Value = -3572.6
$ws.Range("H86").Value = 1572.25
$ws.Range("I86").Value = 1572.25
$ws.Range("K86").Value = 1572.25
$ws.Range("M86").Value = -449.25
$ws.Range("H89").Value = 1572.25
$ws.Range("I89").Value = 1572.25
$ws.Range("K89").Value = 7861.25
$ws.Range("M89").Value = -2245.25
$ws.Range("H94").Value = 1271.75
$ws.Range("I94").Value = 1179.6666
$ws.Range("K94").Value = 1179.6666
$ws.Range("M94").Value = -728.6666
$ws.Range("H99").Value = 1933.1666
$ws.Range("I99").Value = 1927.091
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 1927.091
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = -429.0909999999999
$ws.Range("N99").Value = -4996
$ws.Range("H106").Value = 7666.3335
$ws.Range("J106").Value = 7666.3335
$ws.Range("L106").Value = 7666.3335
$ws.Range("N106").Value = -10190.3335
$ws.Range("H134").Value = 1882.9166
$ws.Range("I134").Value = 1834.3478
$ws.Range("K134").Value = 5503.0434
$ws.Range("M134").Value = -2968.0434

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2581.125
$ws.Range("I58").Value = 1996.8
$ws.Range("J58").Value = 3555
$ws.Range("K58").Value = 1996.8
$ws.Range("L58").Value = 3555
$ws.Range("M58").Value = -1793.8
$ws.Range("N58").Value = -3961
$ws.Range("H99").Value = 1899
$ws.Range("I99").Value = 1771.8182
$ws.Range("K99").Value = 1771.8182
$ws.Range("M99").Value = -273.8181999999999
$ws.Range("H126").Value = 1899
$ws.Range("I126").Value = 1771.8182
$ws.Range("K126").Value = 5315.4546
$ws.Range("M126").Value = -2845.4546
$ws.Range("H134").Value = 2293.5881
$ws.Range("I134").Value = 2298.0625
$ws.Range("J134").Value = 2222
$ws.Range("K134").Value = 6894.1875
$ws.Range("L134").Value = 6666
$ws.Range("M134").Value = -4359.1875
$ws.Range("N134").Value = -11736
$ws.Range("H136").Value = 2581.125
$ws.Range("I136").Value = 1996.8
$ws.Range("J136").Value = 3555
$ws.Range("K136").Value = 5990.4
$ws.Range("L136").Value = 10665
$ws.Range("M136").Value = -3440.4
$ws.Range("N136").Value = -15765

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 750
$ws.Range("J86").Value = 750
$ws.Range("L86").Value = 2250
$ws.Range("N86").Value = -4622
$ws.Range("H89").Value = 750
$ws.Range("J89").Value = 750
$ws.Range("L89").Value = 6750
$ws.Range("N89").Value = -18606
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").ClearContents()
$ws.Range("N95").Value = 0
$ws.Range("H113").Value = 1256.6154
$ws.Range("J113").Value = 1368.375
$ws.Range("L113").Value = 4105.125
$ws.Range("N113").Value = -8445.125
$ws.Range("H131").Value = 2588.7778

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H94").Value = 31159.1
$ws.Range("J94").Value = 31159.1
$ws.Range("L94").Value = 31159.1
$ws.Range("N94").Value = -32511.1
$ws.Range("H105").Value = 23390.5
$ws.Range("J105").Value = 23390.5
$ws.Range("L105").Value = 23390.5
$ws.Range("N105").Value = -30378.5
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H132").Value = 10927.846
$ws.Range("I132").Value = 10927.846
$ws.Range("K132").Value = 32783.538
$ws.Range("M132").Value = -30253.538

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 706.7143
$ws.Range("J22").Value = 399.5
$ws.Range("L22").Value = 399.5
$ws.Range("N22").Value = -989.5
$ws.Range("H27").Value = 706.7143
$ws.Range("J27").Value = 399.5
$ws.Range("L27").Value = 399.5
$ws.Range("N27").Value = -613.5
$ws.Range("H40").Value = 4674.7
$ws.Range("I40").Value = 4178.143
$ws.Range("K40").Value = 4178.143
$ws.Range("M40").Value = -4042.143
$ws.Range("H98").Value = 22500
$ws.Range("J98").Value = 22500
$ws.Range("L98").Value = 22500
$ws.Range("N98").Value = -28490
$ws.Range("H105").Value = 28750
$ws.Range("J105").Value = 28750
$ws.Range("L105").Value = 28750
$ws.Range("N105").Value = -35738
$ws.Range("H138").Value = 90000
$ws.Range("J138").Value = 90000
$ws.Range("L138").Value = 90000
$ws.Range("N138").Value = -100280

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 2250
$ws.Range("I9").Value = 4006
$ws.Range("K9").Value = 4006
$ws.Range("M9").Value = -3866
$ws.Range("H62").Value = 17000.428
$ws.Range("J62").Value = 18666.334
$ws.Range("L62").Value = 18666.334
$ws.Range("N62").Value = -19914.334
$ws.Range("H65").Value = 17000.428
$ws.Range("J65").Value = 18666.334
$ws.Range("L65").Value = 93331.67
$ws.Range("N65").Value = -99571.67
$ws.Range("H107").Value = 866.3333
$ws.Range("I107").Value = 849
$ws.Range("K107").Value = 2547
$ws.Range("M107").Value = -627
